# Add Test Data for Russia, Finland and Hungary markets.
# Modeled after the "Croatia" sheet, which already has a matching layout
# (empty B3 / filled B4, column A/B widths matching the new markets' sheets).

$wb = $excel.ActiveWorkbook

function Add-MarketSheet {
    param($TemplateName, $NewName, $CodeValue, $MarketValue)

    $template = $wb.Worksheets.Item($TemplateName)
    [void]$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $NewName

    # Set B4 (code) before B2 (market name) so the new shared strings are
    # appended in the same order the source workbook used.
    $newSheet.Range("B4").Value = $CodeValue
    $newSheet.Range("B2").Value = $MarketValue

    # Match the row heights used by the new sheets (wrapped 2-line rows).
    $newSheet.Rows.Item(3).RowHeight = 28.8
    $newSheet.Rows.Item(4).RowHeight = 28.8
    $newSheet.Rows.Item(5).RowHeight = 28.8

    return $newSheet
}

$russia = Add-MarketSheet "Croatia" "Russia" "NGC-2929/T2900 " "Russia Market"
[void]$russia.Range("A1:D11").Select()

$finland = Add-MarketSheet "Croatia" "Finland" "NGC-3130/T2943 " "Finland Market"
[void]$finland.Range("A1:D11").Select()

$hungary = Add-MarketSheet "Croatia" "Hungary" "NGC-3104/T2992 " "Hungary Market"
[void]$hungary.Range("J17").Select()
